$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the style used by the other
# header cells (e.g. G1: bold, bordered, centered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column data
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
